$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Hoja2"

# Column A - municipi (string) in an order that reproduces the target sharedStrings sequence
$ws2.Cells.Item(1, 1).Value = "municipi"
$ws2.Cells.Item(2, 1).Value = "Abrera"
$ws2.Cells.Item(3, 1).Value = "Aguilar"
$ws2.Cells.Item(4, 1).Value = "Alélla"
$ws2.Cells.Item(5, 1).Value = "Alpans"
$ws2.Cells.Item(6, 1).Value = "Ametlla Vallès"
$ws2.Cells.Item(8, 1).Value = "L'espunyola"
$ws2.Cells.Item(7, 1).Value = "Arenys de Mar"
$ws2.Cells.Item(9, 1).Value = " Terrassa "
$ws2.Cells.Item(10, 1).Value = "Barcelona "
$ws2.Cells.Item(11, 1).Value = "Barcelona "
$ws2.Cells.Item(12, 1).Value = "Alpans"
$ws2.Cells.Item(13, 1).Value = "Alpans"
$ws2.Cells.Item(14, 1).Value = "Alpans13"

# Column B - info (header string, rest numeric)
$ws2.Cells.Item(1, 2).Value = "info"
$ws2.Cells.Item(2, 2).Value = 1
$ws2.Cells.Item(3, 2).Value = 2
$ws2.Cells.Item(4, 2).Value = 3
$ws2.Cells.Item(5, 2).Value = 4
$ws2.Cells.Item(6, 2).Value = 5
$ws2.Cells.Item(7, 2).Value = 6
$ws2.Cells.Item(8, 2).Value = 7
$ws2.Cells.Item(9, 2).Value = 8
$ws2.Cells.Item(10, 2).Value = 9
$ws2.Cells.Item(11, 2).Value = 10
$ws2.Cells.Item(12, 2).Value = 11
$ws2.Cells.Item(13, 2).Value = 12
$ws2.Cells.Item(14, 2).Value = 13

# View state: select A1:B10 on sheet1 (previously active, now not) and H12 on sheet2 (now active)
$ws1.Range("A1:B10").Select()
$ws2.Range("H12").Select()

Write-Host "done"
